$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths (target stored widths: B=43.42578125, C=25.5703125, D=52.28515625).
# ColumnWidth is rounded internally to the nearest 1/6 of a character (pixel grid),
# so these inputs are chosen to land on the closest representable stored width.
$ws.Columns.Item(2).ColumnWidth = 42.66666666666667
$ws.Columns.Item(3).ColumnWidth = 24.666666666666664
$ws.Columns.Item(4).ColumnWidth = 51.5

# Column A (Code)
$ws.Range("A2").Value = "GINF1"
$ws.Range("A3").Value = "GINF52"
$ws.Range("A4").Value = "GINF53"
$ws.Range("A5").Value = "GINF54"
$ws.Range("A6").Value = "GINF55"
$ws.Range("A7").Value = "GINF56"

# Column B (Intitulé)
$ws.Range("B2").Value = "Technologies Net et J2EE"
$ws.Range("B3").Value = "IA Avancée et Ingénierie de connaissance"
$ws.Range("B4").Value = "Systèmes d'information et BI"
$ws.Range("B5").Value = "Management des SI"
$ws.Range("B6").Value = "Web Services et applications"
$ws.Range("B7").Value = "Management de l'entreprise"

# Column C (Chef Module)
$ws.Range("C2").Value = "El Haddad"
$ws.Range("C3").Value = "El Alami Hassoun"
$ws.Range("C4").Value = "Badir"
$ws.Range("C5").Value = "Ezzine"
$ws.Range("C6").Value = "Ben Achhab"
$ws.Range("C7").Value = "El Haddad"

# Column D (Composants) - written in this specific order to match original shared-string ordering
$ws.Range("D2").Value = "Architecture J2EE, Prog Mobile Net"
$ws.Range("D4").Value = "Datawarehouse, Dtatamining, BigData et Applications"
$ws.Range("D5").Value = "Gouvernance des SI, ERP & CRM, Audit & Sécurité des SI"
$ws.Range("D3").Value = "IA Avancée & Web Sémantique, Ingénierie de connaissances"
$ws.Range("D6").Value = "Dév Web en .net, Webservices, E-commerce et applications"
$ws.Range("D7").Value = "Création d'Entreprise, Projets libres, Simulation d'entretien d'embauche et éthique de l'ingénieur"

# Selection
$ws.Range("C7").Select()
